$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Holly added "S.GISH" as a harvester in bioSamples, so the rnaSamples
# sheet's "harvester" column (column B, rows 2-25) needs to be fixed to
# reference that same harvester value instead of the old placeholder
# text that had been sitting there ("Retrofitted_0629_0618").
$ws.Range("B2:B25").Value = "S.GISH"

# Reflect that column B was the column being worked on / selected, and
# nudge its width slightly (as happens when a column's content is
# edited/autofit in Excel).
$ws.Columns("B").ColumnWidth = 7.9967
$ws.Columns("B").Select()
